$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("I2").Value = 8
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 2.5
$ws.Range("W2").Value = 9
$ws.Range("Z2").Value = 9.5
$ws.Range("AB2").Value = 21
$ws.Range("AK2").Value = 81

# Row 4
$ws.Range("G4").Value = 1.85
$ws.Range("I4").Value = 5.25
$ws.Range("J4").Value = 2.6
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("X4").Value = 7
$ws.Range("AH4").Value = 10
$ws.Range("AI4").Value = 23
$ws.Range("AJ4").Value = 19
$ws.Range("AL4").Value = 51
$ws.Range("AU4").Value = 10
$ws.Range("AW4").Value = 6.5
$ws.Range("AZ4").Value = 126

# Row 5
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 1.5

# Row 6
$ws.Range("H6").Value = 2.7
$ws.Range("I6").Value = 2.65

# Row 7
$ws.Range("G7").Value = 3.2
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 2.35
$ws.Range("J7").Value = 4
$ws.Range("L7").Value = 3.25
$ws.Range("W7").Value = 7.5
$ws.Range("X7").Value = 13
$ws.Range("Y7").Value = 12
$ws.Range("Z7").Value = 34
$ws.Range("AI7").Value = 10
$ws.Range("AJ7").Value = 10
$ws.Range("AK7").Value = 23
$ws.Range("AL7").Value = 23
$ws.Range("AN7").Value = 5
$ws.Range("AO7").Value = 19
$ws.Range("AQ7").Value = 67
$ws.Range("AW7").Value = 4.33
$ws.Range("BB7").Value = 251

# Row 8
$ws.Range("G8").Value = 1.71
$ws.Range("I8").Value = 5
$ws.Range("AF8").Value = 81
$ws.Range("AJ8").Value = 17
$ws.Range("AX8").Value = 29
$ws.Range("BA8").Value = 151

# Row 10
$ws.Range("K10").Value = 2.37
